$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_of_studies")

# --- Add a new study record in row 41 ---
$ws.Range("A41").Value = "The Impact of Virtual Reality on Chronic Pain"
$ws.Range("B41").Value = 2016
$ws.Range("C41").Value = "PloS One"
$ws.Range("D41").Value = "Jones, Ted; Moore, Todd; Choo, James"
$ws.Range("F41").Value = "Pain"
$ws.Range("G41").Value = "Patient Study"
$ws.Range("K41").Value = "Experiment"

# Link column (E41): add the hyperlink, then restore the standard
# "Hyperlink" cell style used throughout the rest of the sheet.
$ws.Hyperlinks.Add($ws.Range("E41"), "https://doi.org/10.1371/journal.pone.0167523") | Out-Null
$ws.Range("E41").Style = "Hyperlink"

# --- Update the sheet view (scroll position + active selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A41").Select() | Out-Null
